# Insert a new weekly record for "Orégano" at Mercado Mayorista Lo Valledor de
# Santiago. This pushes all existing observations starting at row 29 down by
# one row (dimension grows from A1:R92 to A1:R93), and the new row 29 is
# populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 29:92 down to 30:93, leaving a blank row 29 to fill in.
$ws.Rows.Item(29).Insert()

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44414
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 100112029
$ws.Range("G29").Value = "Orégano"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 29
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 9483
$ws.Range("N29").Value = "$/docena de atados"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 3161
$ws.Range("Q29").Value = 3
$ws.Range("R29").Value = "Hortaliza"
